$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H112").Value = 1486.9166
$ws.Range("J112").Value = 1486.9166
$ws.Range("L112").Value = 4460.7498
$ws.Range("N112").Value = -6676.7498
$ws.Range("H132").Value = 812.1429000000001
$ws.Range("I132").Value = 693.5517
$ws.Range("K132").Value = 2080.6551
$ws.Range("M132").Value = 449.3449000000001

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H2").Value = 293484.3
$ws.Range("I2").Value = 397895
$ws.Range("J2").Value = 1134.4
$ws.Range("K2").Value = 397895
$ws.Range("L2").Value = 1134.4
$ws.Range("M2").Value = -397782
$ws.Range("N2").Value = -1360.4
$ws.Range("H45").Value = 5001351.5
$ws.Range("I45").Value = 11250843
$ws.Range("J45").Value = 1758.3
$ws.Range("K45").Value = 11250843
$ws.Range("L45").Value = 1758.3
$ws.Range("M45").Value = -11250466
$ws.Range("N45").Value = -2512.3
$ws.Range("H74").Value = 1113.6666
$ws.Range("I74").Value = 500
$ws.Range("K74").Value = 500
$ws.Range("M74").Value = 374
$ws.Range("H77").Value = 1113.6666
$ws.Range("I77").Value = 500
$ws.Range("K77").Value = 2500
$ws.Range("M77").Value = 1868
$ws.Range("H116").Value = 293484.3
$ws.Range("I116").Value = 397895
$ws.Range("J116").Value = 1134.4
$ws.Range("K116").Value = 397895
$ws.Range("L116").Value = 1134.4
$ws.Range("M116").Value = -395601
$ws.Range("N116").Value = -5722.4
$ws.Range("H132").Value = 1649.1875
$ws.Range("I132").Value = 1302.625
$ws.Range("K132").Value = 3907.875
$ws.Range("M132").Value = -1377.875

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H3").Value = 293484.3
$ws.Range("I3").Value = 397895
$ws.Range("J3").Value = 1134.4
$ws.Range("K3").Value = 397895
$ws.Range("L3").Value = 1134.4
$ws.Range("M3").Value = -397781
$ws.Range("N3").Value = -1362.4
$ws.Range("H12").Value = 16249.25
$ws.Range("I12").Value = 11666
$ws.Range("J12").Value = 29999
$ws.Range("K12").Value = 11666
$ws.Range("L12").Value = 29999
$ws.Range("M12").Value = -11498
$ws.Range("N12").Value = -30335
$ws.Range("H99").Value = 1775
$ws.Range("I99").Value = 1109
$ws.Range("J99").Value = 1997
$ws.Range("K99").Value = 1109
$ws.Range("L99").Value = 1997
$ws.Range("M99").Value = 389
$ws.Range("N99").Value = -4993
$ws.Range("H105").Value = 2588.762
$ws.Range("I105").Value = 2242.4443
$ws.Range("K105").Value = 2242.4443
$ws.Range("M105").Value = -495.4443000000001
$ws.Range("H134").Value = 5182.3447
$ws.Range("I134").Value = 5598.5415
$ws.Range("J134").Value = 3184.6
$ws.Range("K134").Value = 16795.6245
$ws.Range("L134").Value = 9553.799999999999
$ws.Range("M134").Value = -14260.6245
$ws.Range("N134").Value = -14623.8

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H16").Value = 691.4545000000001
$ws.Range("J16").Value = 600.25
$ws.Range("L16").Value = 600.25
$ws.Range("N16").Value = -1174.25
$ws.Range("H31").Value = 2400
$ws.Range("I31").Value = 4000
$ws.Range("K31").Value = 4000
$ws.Range("M31").Value = -3705
$ws.Range("H34").Value = 2400
$ws.Range("I34").Value = 4000
$ws.Range("K34").Value = 4000
$ws.Range("M34").Value = -3798
$ws.Range("H58").Value = 4833462
$ws.Range("I58").Value = 10870240
$ws.Range("K58").Value = 10870240
$ws.Range("M58").Value = -10870037
$ws.Range("H59").Value = 18100
$ws.Range("J59").Value = 18100
$ws.Range("L59").Value = 18100
$ws.Range("N59").Value = -20390
$ws.Range("H107").Value = 915.13336
$ws.Range("J107").Value = 1439.8
$ws.Range("L107").Value = 1439.8
$ws.Range("N107").Value = -5279.8
$ws.Range("H113").Value = 691.4545000000001
$ws.Range("J113").Value = 600.25
$ws.Range("L113").Value = 600.25
$ws.Range("N113").Value = -4940.25
$ws.Range("H132").Value = 2491.182
$ws.Range("I132").Value = 1519.8667
$ws.Range("K132").Value = 4559.6001
$ws.Range("M132").Value = -2029.6001
$ws.Range("H136").Value = 4833462
$ws.Range("I136").Value = 10870240
$ws.Range("K136").Value = 32610720
$ws.Range("M136").Value = -32608170
$ws.Range("H140").Value = 0
$ws.Range("J140").Value = 0
$ws.Range("L140").Value = 0
$ws.Range("N140").ClearContents()

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H131").Value = 12172.528
$ws.Range("J131").Value = 12508.044
$ws.Range("L131").Value = 37524.132
$ws.Range("N131").Value = -47604.132

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H113").Value = 1434.375
$ws.Range("I113").Value = 1151
$ws.Range("K113").Value = 1151
$ws.Range("M113").Value = 1019

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H7").Value = 2423.5557
$ws.Range("I7").Value = 2320.25
$ws.Range("K7").Value = 2320.25
$ws.Range("M7").Value = -2208.25
$ws.Range("H100").Value = 1677.1666
$ws.Range("I100").Value = 894.3333
$ws.Range("J100").Value = 2460
$ws.Range("K100").Value = 894.3333
$ws.Range("M100").Value = -353.3333
$ws.Range("N100").Value = -3542
$ws.Range("H126").Value = 2423.5557
$ws.Range("I126").Value = 2320.25
$ws.Range("K126").Value = 6960.75
$ws.Range("M126").Value = -4490.75
$ws.Range("H132").Value = 1891.814
$ws.Range("J132").Value = 2309.32
$ws.Range("L132").Value = 6927.960000000001
$ws.Range("N132").Value = -11987.96
$ws.Range("H136").Value = 4987.5
$ws.Range("I136").Value = 3379
$ws.Range("J136").Value = 7668.3335
$ws.Range("K136").Value = 10137
$ws.Range("L136").Value = 23005.0005
$ws.Range("M136").Value = -7587
$ws.Range("N136").Value = -28105.0005

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H107").Value = 811
$ws.Range("I107").Value = 573.2
$ws.Range("J107").Value = 2000
$ws.Range("K107").Value = 1719.6
$ws.Range("L107").Value = 6000
$ws.Range("M107").Value = 200.3999999999999
$ws.Range("N107").Value = -9840
$ws.Range("H113").Value = 524.7778
$ws.Range("I113").Value = 354.35715
$ws.Range("K113").Value = 1063.07145
$ws.Range("M113").Value = 1106.92855
$ws.Range("H122").Value = 28589.072
$ws.Range("I122").Value = 32987.543
$ws.Range("K122").Value = 98962.62899999999
$ws.Range("M122").Value = -96512.62899999999
$ws.Range("H132").Value = 1047.619
$ws.Range("I132").Value = 815.7083
$ws.Range("J132").Value = 1789.7333
$ws.Range("K132").Value = 2447.1249
$ws.Range("L132").Value = 5369.199900000001
$ws.Range("M132").Value = 82.8751000000002
$ws.Range("N132").Value = -10429.1999
